$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "power_line_Wholesale_Kasso" row (row 8), shifting subsequent rows up
$ws.Rows.Item(8).Delete()
